# "corrected incorrect lists ;)"
# The workbook had a number of mixed-up shared-string entries in columns
# B:G of worksheet "Tabelle1" (swapped between cells in the same row).
# Fix them so each row reads correctly again, and restore the view state
# (active sheet/selection) that was in effect when the correction was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# --- Row 4 ---
$ws.Range("C4").Value = "ki"
$ws.Range("D4").Value = "we"
$ws.Range("E4").Value = "ki"
$ws.Range("G4").Value = "lu"

# --- Row 12 ---
$ws.Range("B12").Value = "lu"
$ws.Range("F12").Value = "pu"

# --- Row 13 ---
$ws.Range("C13").Value = "lu"
$ws.Range("G13").Value = "ki"

# --- Row 17 ---
$ws.Range("C17").Value = "we"
$ws.Range("D17").Value = "we"
$ws.Range("E17").Value = "me"

# --- Row 20 ---
$ws.Range("F20").Value = "lu"
$ws.Range("G20").Value = "ki"

# --- Row 21 ---
$ws.Range("C21").Value = "go"
$ws.Range("D21").Value = "lu"

# --- Row 22 ---
$ws.Range("D22").Value = "fo"
$ws.Range("F22").Value = "lu"

# --- Row 25 ---
$ws.Range("C25").Value = "si"
$ws.Range("F25").Value = "ta"

# --- Row 27 ---
$ws.Range("E27").Value = "lu"
$ws.Range("F27").Value = "ki"

# --- Row 28 ---
$ws.Range("B28").Value = "me"
$ws.Range("C28").Value = "na"

# --- Row 32 ---
$ws.Range("E32").Value = "ta"
$ws.Range("F32").Value = "ki"

# --- Row 33 ---
$ws.Range("F33").Value = "ta"
$ws.Range("G33").Value = "me"

# --- Row 34 ---
$ws.Range("D34").Value = "we"
$ws.Range("F34").Value = "ta"
$ws.Range("G34").Value = "lu"

# --- Row 40 ---
$ws.Range("D40").Value = "lu"
$ws.Range("F40").Value = "si"

# Restore the on-screen selection: the sheet was scrolled down and the
# whole of column A was selected.
$ws.Columns.Item(1).Select()
